{"js": "// Update the worksheet answer cells (three-digit \u00f7 one-digit division problems)\n// to the new set of problems/answers, per cell position. Using direct\n// table-cell addressing (row, col) keeps each edit targeted at the exact\n// cell regardless of whether any text value is duplicated elsewhere in the\n// table (some of the new values coincide with old values from a different\n// cell).\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Each entry: [rowIndex, colIndex, newText]\nconst updates = [\n  [0, 0, \"225\u00f75=45, 0\"],\n  [0, 1, \"658\u00f77=94, 0\"],\n  [0, 2, \"186\u00f78=23, 2\"],\n  [0, 3, \"642\u00f76=107, 0\"],\n  [0, 4, \"878\u00f76=146, 2\"],\n\n  [4, 0, \"422\u00f74=105, 2\"],\n  [4, 1, \"528\u00f72=264, 0\"],\n  [4, 2, \"696\u00f75=139, 1\"],\n  [4, 3, \"298\u00f74=74, 2\"],\n  [4, 4, \"663\u00f77=94, 5\"],\n\n  [8, 0, \"460\u00f75=92, 0\"],\n  [8, 1, \"907\u00f79=100, 7\"],\n  [8, 2, \"114\u00f72=57, 0\"],\n  [8, 3, \"763\u00f79=84, 7\"],\n  [8, 4, \"261\u00f73=87, 0\"],\n\n  [12, 0, \"293\u00f72=146, 1\"],\n  [12, 1, \"769\u00f75=153, 4\"],\n  [12, 2, \"292\u00f74=73, 0\"],\n  [12, 3, \"128\u00f73=42, 2\"],\n  [12, 4, \"442\u00f72=221, 0\"],\n\n  [16, 0, \"759\u00f72=379, 1\"],\n  [16, 1, \"848\u00f75=169, 3\"],\n  [16, 2, \"793\u00f78=99, 1\"],\n  [16, 3, \"165\u00f74=41, 1\"],\n  [16, 4, \"845\u00f73=281, 2\"],\n];\n\nfor (const [row, col, text] of updates) {\n  table.getCell(row, col).value = text;\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet answer cells (three-digit \u00f7 one-digit division\n# problems) to the new set of problems/answers, per cell position. Using\n# direct table-cell addressing (row, col) keeps each edit targeted at the\n# exact cell regardless of whether any text value is duplicated elsewhere\n# in the table (some of the new values coincide with old values from a\n# different cell).\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Each entry: row, column (1-based, matching Word's Table.Cell(r, c)), newText\n$updates = @(\n    @(1, 1, \"225\u00f75=45, 0\"),\n    @(1, 2, \"658\u00f77=94, 0\"),\n    @(1, 3, \"186\u00f78=23, 2\"),\n    @(1, 4, \"642\u00f76=107, 0\"),\n    @(1, 5, \"878\u00f76=146, 2\"),\n\n    @(5, 1, \"422\u00f74=105, 2\"),\n    @(5, 2, \"528\u00f72=264, 0\"),\n    @(5, 3, \"696\u00f75=139, 1\"),\n    @(5, 4, \"298\u00f74=74, 2\"),\n    @(5, 5, \"663\u00f77=94, 5\"),\n\n    @(9, 1, \"460\u00f75=92, 0\"),\n    @(9, 2, \"907\u00f79=100, 7\"),\n    @(9, 3, \"114\u00f72=57, 0\"),\n    @(9, 4, \"763\u00f79=84, 7\"),\n    @(9, 5, \"261\u00f73=87, 0\"),\n\n    @(13, 1, \"293\u00f72=146, 1\"),\n    @(13, 2, \"769\u00f75=153, 4\"),\n    @(13, 3, \"292\u00f74=73, 0\"),\n    @(13, 4, \"128\u00f73=42, 2\"),\n    @(13, 5, \"442\u00f72=221, 0\"),\n\n    @(17, 1, \"759\u00f72=379, 1\"),\n    @(17, 2, \"848\u00f75=169, 3\"),\n    @(17, 3, \"793\u00f78=99, 1\"),\n    @(17, 4, \"165\u00f74=41, 1\"),\n    @(17, 5, \"845\u00f73=281, 2\")\n)\n\nforeach ($u in $updates) {\n    $row = $u[0]\n    $col = $u[1]\n    $text = $u[2]\n    $t.Cell($row, $col).Range.Text = $text\n}\n"}
